# feat: Add modern KPI dashboard and enhance terminal output
#
# Renames the original sheet to "Sales Data" and appends a new
# "Summary Statistics" sheet holding descriptive statistics
# (pandas-style .describe()) computed over the "Sales" column.

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet -------------------------------------------------
$salesSheet = $wb.Worksheets.Item(1)
$salesSheet.Name = "Sales Data"

# --- Add the new summary sheet right after the sales data sheet ---------------
$summarySheet = $wb.Worksheets.Add($null, $salesSheet)
$summarySheet.Name = "Summary Statistics"

# --- Pre-format the label column as text so values such as "25%" are kept as
#     literal text instead of being auto-parsed into numeric percentages. ------
$labelRange = $summarySheet.Range("A2:A9")
$labelRange.NumberFormat = "@"

# --- Header -------------------------------------------------------------------
$summarySheet.Cells.Item(1, 2).Value = "Sales"

# --- Descriptive statistics rows ---------------------------------------------
$stats = @(
    @("count", 4),
    @("mean", 2875),
    @("std", 853.9125638299665),
    @("min", 2000),
    @("25%", 2375),
    @("50%", 2750),
    @("75%", 3250),
    @("max", 4000)
)

$row = 2
foreach ($pair in $stats) {
    $summarySheet.Cells.Item($row, 1).Value = $pair[0]
    $summarySheet.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# --- Styling: reuse the bold / centered / thin-bordered header style already
#     defined on "Sales Data"!A1 so no new styles are introduced. -------------
$salesSheet.Range("A1").Copy()
$labelRange.PasteSpecial(-4122)

$salesSheet.Range("B1").Copy()
$summarySheet.Range("B1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$summarySheet.Range("A1").Select() | Out-Null
